# Daily attendance processing - 2025-12-02 19:27:31
# Normalizes the "Recorded By" (column G) values on the session-analysis
# sheet by alphabetically (ordinal) sorting the comma-separated list of
# recorder names/emails in each row.

function CompareOrdinal($a, $b) {
    $la = $a.Length
    $lb = $b.Length
    $m = $la
    if ($lb -lt $m) { $m = $lb }
    $i = 0
    while ($i -lt $m) {
        $ca = [int][char]$a[$i]
        $cb = [int][char]$b[$i]
        if ($ca -ne $cb) { return $ca - $cb }
        $i = $i + 1
    }
    return $la - $lb
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $t = $ws.Cells.Item($r, 7).Text
    if ([string]::IsNullOrEmpty($t)) {
        continue
    }

    $parts = $t -split ", "
    $n = $parts.Count

    if ($n -eq 2) {
        if ((CompareOrdinal $parts[0] $parts[1]) -gt 0) {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
        }
    } elseif ($n -eq 3) {
        if ((CompareOrdinal $parts[0] $parts[1]) -gt 0) {
            $tmp = $parts[0]; $parts[0] = $parts[1]; $parts[1] = $tmp
        }
        if ((CompareOrdinal $parts[1] $parts[2]) -gt 0) {
            $tmp = $parts[1]; $parts[1] = $parts[2]; $parts[2] = $tmp
        }
        if ((CompareOrdinal $parts[0] $parts[1]) -gt 0) {
            $tmp = $parts[0]; $parts[0] = $parts[1]; $parts[1] = $tmp
        }
    }

    if ($n -ge 2) {
        $joined = [string]::Join(", ", $parts)
        if ($joined -ne $t) {
            $ws.Cells.Item($r, 7).Value = $joined
        }
    }
}
